$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 89
$ws.Range("I2").Value = 196
$ws.Range("J2").Value = 871
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 216
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 145
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 92
$ws.Range("T2").Value = 132
$ws.Range("U2").Value = 10
$ws.Range("V2").Value = 1215
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 1286
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 16
$ws.Range("AA2").Value = 11
